$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Subtitle byline: collapse the spell-check-wrapped runs into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("18740: Shravani Dhote, Simrit Kaur, Vins Sharma", $false, $false, $false, $false, $false, $true, 1, $false, "18740: Shravani Dhote, Simrit Kaur, Vins Sharma", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "associativity, block sizing, set sizing, etc, ..." run merge.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("associativity, block sizing, set sizing, etc, and still reuse the base code, we cut down the number of ways to ", $false, $false, $false, $false, $false, $true, 1, $false, "associativity, block sizing, set sizing, etc, and still reuse the base code, we cut down the number of ways to ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) ", without changing the number of sets, block sizings, ..." run merge.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(", without changing the number of sets, block sizings, or any significant functionality.", $false, $false, $false, $false, $false, $true, 1, $false, ", without changing the number of sets, block sizings, or any significant functionality.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) New paragraph introducing the scheduler graphs, right before the
#    "Equity scheduler" Heading2 paragraph.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Task 5: A custom scheduler", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor4 = $r4.Paragraphs.First
$anchor4.Range.InsertParagraphAfter()
$newPara4 = $d.Paragraphs.Item($anchor4.Index + 1)
$xml4 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The following are a series of graphs of </w:t></w:r><w:r><w:t>all of our scheduling mechanisms, overlaid upon each other to compare performance.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 5) Three new paragraphs appended at the end of the document: the
#    "This resulted in..." wrap-up, the "Equitable BLISS scheduler" Heading2,
#    and its description paragraph.
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Collapse(0)
$xml5 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This resulted in a decent amount of performance, </w:t></w:r><w:r><w:t>but it’s not as fair as way partitioning.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Equitable BLISS scheduler</w:t></w:r></w:p><w:p><w:r><w:t>Our next idea was to leverage</w:t></w:r><w:r><w:t xml:space="preserve"> BLISS’s performance gains but still add another key of further incentivizing low usage to particular cores. As such, we implemented BLISS alongside our equity scheduler, with our equity scheduler being the tiebreaker for any BLISS operations and having FR-FCFS as our final tiebreaker. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r5.InsertXML($xml5)
